# Shift each schedule date in column B (rows 2-20) forward by 9 days.
# Dates are stored as plain text strings in "dd.MM.yyyy" format, so we must
# force text entry (some of the new values, e.g. "01.12.2024", would
# otherwise be auto-recognized by Excel as real dates).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$newDates = @(
    "15.11.2024",
    "16.11.2024",
    "17.11.2024",
    "18.11.2024",
    "19.11.2024",
    "20.11.2024",
    "21.11.2024",
    "22.11.2024",
    "23.11.2024",
    "24.11.2024",
    "25.11.2024",
    "26.11.2024",
    "27.11.2024",
    "28.11.2024",
    "29.11.2024",
    "30.11.2024",
    "01.12.2024",
    "02.12.2024",
    "03.12.2024"
)

for ($i = 0; $i -lt $newDates.Length; $i++) {
    $row = $i + 2
    $cell = $ws.Cells.Item($row, 2)
    $cell.NumberFormat = "@"
    $cell.Value = $newDates[$i]
    $cell.Style = "Normal"
}
